# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# as scraped/refreshed by the GitHub Actions workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.034.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.453.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '510.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("E8").Value = '  -1.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.453.52'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0980'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.20%  '

$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("E13").Value = '  -7.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.889.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.886.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("E17").Value = '  +1.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.403.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '314.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.70%  '

$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = '  -1.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.380'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '172.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0732'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.13'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.26%  '

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.03%  '

$ws.Range("E38").Value = '  +5.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.70'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.47'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.807'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.94'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.576'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '256.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.65%  '

$ws.Range("E48").Value = '  -0.55%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0493'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0214'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.45%  '
